$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the Subtotal row (old row 8) to host the two new
# budget line items (rows 8 and 9 in the new layout).
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# Copy the formatting from an existing data row (row 4) into the two new rows
# so the new cells pick up the same fonts/borders/number formats.
$ws.Range("A4:G4").Copy()
$ws.Range("A8:G9").PasteSpecial(-4122)

$h4 = $ws.Rows.Item(4).RowHeight
$ws.Rows.Item(8).RowHeight = $h4
$ws.Rows.Item(9).RowHeight = $h4

# New row 8: 6-pin Headers
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "6-pin Headers"
$ws.Range("C8").Value = "Creatron Inc"
$ws.Range("E8").Value = "To solder sensor into PCB"
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 1.78

# New row 9: 2x20 pin Headers
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "2x20 pin Headers"
$ws.Range("C9").Value = "Creatron Inc"
$ws.Range("E9").Value = "For Raspberry Pi"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.89

# Part numbers filled in last (matches authoring order).
$ws.Range("D8").Value = "CONHD-330060"
$ws.Range("D9").Value = "CONPH-402749"

# Update the Subtotal / HST / Total formulas (now on rows 10, 11, 12) so they
# take the two new line items into account.
$ws.Range("G10").Formula = "=SUM(G4:G9)"
$ws.Range("G11").Formula = "=G10*0.13"
$ws.Range("G12").Formula = "=SUM(G10:G11)"

# Widen column E to fit the longer description text.
$ws.Columns.Item(5).ColumnWidth = 30.33203125

# Match the selection left behind by the author after editing.
$ws.Range("G10").Select()
